# Insert a new data row at row 344 (shifting existing rows 344-427 down to 345-428)
# and populate it with a new Albahaca price record for "Región Metropolitana".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(344).Insert()

$ws.Range("A344").Value2 = 9
$ws.Range("B344").Value2 = "Vega Central Mapocho de Santiago"
$ws.Range("C344").Value2 = "Metropolitana"
$ws.Range("D344").Value2 = 44889
$ws.Range("E344").Value2 = 13
$ws.Range("F344").Value2 = 100112052
$ws.Range("G344").Value2 = "Albahaca"
$ws.Range("H344").Value2 = "Sin especificar"
$ws.Range("I344").Value2 = "Primera"
$ws.Range("J344").Value2 = 310
$ws.Range("K344").Value2 = 5000
$ws.Range("L344").Value2 = 7000
$ws.Range("M344").Value2 = 6161
$ws.Range("N344").Value2 = "`$/docena de matas"
$ws.Range("O344").Value2 = "Región Metropolitana"
$ws.Range("P344").Value2 = 1027
$ws.Range("Q344").Value2 = 6
$ws.Range("R344").Value2 = "Hortaliza"
